$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# account_num: was the repr of a function reference, now a generated numeric
# account number. It must stay a text value (matches the source inlineStr
# type), so prefix with an apostrophe to force text entry instead of Excel
# auto-converting it to a number. Re-apply the Normal style afterwards so we
# don't leave a stray "quote prefix" number format on the cell.
$ws.Range("A2").Value = "'5392379531"
$ws.Range("A2").Style = "Normal"

# username: capitalize
$ws.Range("B2").Value = "Connor"

# password hash: updated hash value
$ws.Range("C2").Value = "ecd71870d1963316a97e3ac3408c9835ad8cf0f3c1bc703527c30265534f75ae"

# date_opened: updated serial date/time value
$ws.Range("F2").Value2 = 46067.57548154639
